# Insert a new data row at row 384 (shifts existing rows 384..410 down to 385..411)
# and populate it with the new weekly price observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(384).Insert()

$ws.Cells(384, 1).Value = 9
$ws.Cells(384, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells(384, 3).Value = "Metropolitana"
$ws.Cells(384, 4).Value = 45075
$ws.Cells(384, 5).Value = 13
$ws.Cells(384, 6).Value = 100112021
$ws.Cells(384, 7).Value = "Ají"
$ws.Cells(384, 8).Value = "Inferno"
$ws.Cells(384, 9).Value = "Primera"
$ws.Cells(384, 10).Value = 34
$ws.Cells(384, 11).Value = 17000
$ws.Cells(384, 12).Value = 19000
$ws.Cells(384, 13).Value = 18000
$ws.Cells(384, 14).Value = "$/caja 10 kilos"
$ws.Cells(384, 15).Value = "Región de Arica y Parinacota"
$ws.Cells(384, 16).Value = 1800
$ws.Cells(384, 17).Value = 10
$ws.Cells(384, 18).Value = "Hortaliza"
